$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K -> E:L
$ws.Columns("D").Insert()

# Copy number formats / styles from column E (the old column D, now shifted) into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest fiscal-year figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1498600
$ws.Range("D9").Value = 1035900
$ws.Range("D10").Value = 462700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 205600
$ws.Range("D17").Value = 1297400
$ws.Range("D18").Value = 201200
$ws.Range("D20").Value = 25300
$ws.Range("D21").Value = 437000
$ws.Range("D22").Value = 67900
$ws.Range("D23").Value = 158600
$ws.Range("D24").Value = -15500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 174100
$ws.Range("D27").Value = 174100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -25300
$ws.Range("D33").Value = 174100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 174100
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 69100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 144400
$ws.Range("D44").Value = 86700
$ws.Range("D45").Value = 34100
$ws.Range("D46").Value = 334300
$ws.Range("D47").Value = 185800
$ws.Range("D48").Value = 3928800
$ws.Range("D49").Value = 223300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 492800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 5165000
$ws.Range("D57").Value = 149800
$ws.Range("D58").Value = 57500
$ws.Range("D59").Value = 197800
$ws.Range("D60").Value = 405100
$ws.Range("D61").Value = 1428500
$ws.Range("D62").Value = 1175600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 3009200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 754600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2155800
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 174100
$ws.Range("D83").Value = 210500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 433100
$ws.Range("D91").Value = -312400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -349000
$ws.Range("D96").Value = -115000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -115200
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -31100
